$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table of updates scraped from the commit diff: row -> (newPrice, newVolume)
# newPrice is $null when the Price (D) column is unchanged for that row.
$updates = @(
    [pscustomobject]@{ Row = 2; Price = '27.786.52'; Volume = '  +1.53%  ' },
    [pscustomobject]@{ Row = 3; Price = '1.885.83'; Volume = '  +1.50%  ' },
    [pscustomobject]@{ Row = 4; Price = '1.008'; Volume = '  +0.60%  ' },
    [pscustomobject]@{ Row = 5; Price = '333.30'; Volume = '  +1.54%  ' },
    [pscustomobject]@{ Row = 6; Price = $null; Volume = '  +0.48%  ' },
    [pscustomobject]@{ Row = 7; Price = '0.4706'; Volume = '  +2.16%  ' },
    [pscustomobject]@{ Row = 8; Price = '0.3937'; Volume = '  +0.04%  ' },
    [pscustomobject]@{ Row = 9; Price = '47.66'; Volume = '  +1.84%  ' },
    [pscustomobject]@{ Row = 10; Price = '0.08076'; Volume = '  +1.83%  ' },
    [pscustomobject]@{ Row = 11; Price = '1.026'; Volume = '  +1.26%  ' },
    [pscustomobject]@{ Row = 12; Price = '22.19'; Volume = '  +3.18%  ' },
    [pscustomobject]@{ Row = 13; Price = '1.891.15'; Volume = '  +1.83%  ' },
    [pscustomobject]@{ Row = 14; Price = '5.985'; Volume = '  +0.88%  ' },
    [pscustomobject]@{ Row = 15; Price = '7.138'; Volume = '  -0.22%  ' },
    [pscustomobject]@{ Row = 16; Price = '1.008'; Volume = '  +0.67%  ' },
    [pscustomobject]@{ Row = 17; Price = '0.06766'; Volume = '  +2.89%  ' },
    [pscustomobject]@{ Row = 18; Price = '87.30'; Volume = '  +0.97%  ' },
    [pscustomobject]@{ Row = 19; Price = '0.00001048'; Volume = '  +1.64%  ' },
    [pscustomobject]@{ Row = 20; Price = $null; Volume = '  +0.93%  ' },
    [pscustomobject]@{ Row = 21; Price = '1.005'; Volume = '  +0.37%  ' },
    [pscustomobject]@{ Row = 22; Price = '27.797.75'; Volume = '  +1.54%  ' },
    [pscustomobject]@{ Row = 23; Price = '5.536'; Volume = '  +0.77%  ' },
    [pscustomobject]@{ Row = 24; Price = '11.03'; Volume = '  +0.93%  ' },
    [pscustomobject]@{ Row = 25; Price = '2.335'; Volume = '  +1.27%  ' },
    [pscustomobject]@{ Row = 26; Price = '2.091.51'; Volume = '  +0.60%  ' },
    [pscustomobject]@{ Row = 27; Price = '159.07'; Volume = '  +3.60%  ' },
    [pscustomobject]@{ Row = 28; Price = '20.16'; Volume = '  +0.35%  ' },
    [pscustomobject]@{ Row = 29; Price = '2.111'; Volume = '  +1.92%  ' },
    [pscustomobject]@{ Row = 30; Price = '5.582'; Volume = '  +1.79%  ' },
    [pscustomobject]@{ Row = 31; Price = '122.03'; Volume = '  +0.32%  ' },
    [pscustomobject]@{ Row = 32; Price = '0.9840'; Volume = '  +2.97%  ' },
    [pscustomobject]@{ Row = 33; Price = '0.09482'; Volume = '  +0.62%  ' },
    [pscustomobject]@{ Row = 34; Price = '1.454'; Volume = '  +0.09%  ' },
    [pscustomobject]@{ Row = 35; Price = '3.618'; Volume = '  +0.80%  ' },
    [pscustomobject]@{ Row = 36; Price = '5.358'; Volume = '  +1.51%  ' },
    [pscustomobject]@{ Row = 37; Price = '0.06148'; Volume = '  +1.76%  ' },
    [pscustomobject]@{ Row = 38; Price = '0.02267'; Volume = '  +1.65%  ' },
    [pscustomobject]@{ Row = 39; Price = '1.226'; Volume = '  +0.43%  ' },
    [pscustomobject]@{ Row = 40; Price = '8.060'; Volume = '  -0.23%  ' },
    [pscustomobject]@{ Row = 41; Price = '0.5999'; Volume = '  +1.16%  ' },
    [pscustomobject]@{ Row = 42; Price = '0.1893'; Volume = '  -0.11%  ' },
    [pscustomobject]@{ Row = 43; Price = '10.30'; Volume = '  +1.32%  ' },
    [pscustomobject]@{ Row = 44; Price = '1.258'; Volume = '  -1.85%  ' },
    [pscustomobject]@{ Row = 45; Price = '0.5720'; Volume = '  +1.40%  ' },
    [pscustomobject]@{ Row = 46; Price = '12.31'; Volume = '  +1.91%  ' },
    [pscustomobject]@{ Row = 47; Price = '1.950'; Volume = '  +1.47%  ' },
    [pscustomobject]@{ Row = 48; Price = '3.396'; Volume = '  -0.03%  ' },
    [pscustomobject]@{ Row = 49; Price = '0.06908'; Volume = '  +2.21%  ' },
    [pscustomobject]@{ Row = 50; Price = '113.66'; Volume = '  +4.55%  ' },
    [pscustomobject]@{ Row = 51; Price = '0.00000000303'; Volume = '  +9.02%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        # Force the Price cell to Text format so numeric-looking strings
        # (e.g. "1.008", "0.00000000303") are kept verbatim instead of being
        # parsed into numbers / scientific notation by Excel.
        $priceCell = $ws.Cells.Item($u.Row, 4)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.Price
    }
    if ($null -ne $u.Volume) {
        $ws.Cells.Item($u.Row, 5).Value = $u.Volume
    }
}
